$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 5688
$ws1.Range("F6").Value = 952
$ws1.Range("F8").Value = 2565
$ws1.Range("F12").Value = 90
$ws1.Range("F13").Value = 30
$ws1.Range("F14").Value = 2406
$ws1.Range("F15").Value = 427

# Sheet "演出" (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 107

# Sheet "全部类型" (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 5688
$ws4.Range("F4").Value = 107
$ws4.Range("F8").Value = 952
$ws4.Range("F10").Value = 2565
$ws4.Range("F15").Value = 90
$ws4.Range("F16").Value = 30
$ws4.Range("F17").Value = 2406
$ws4.Range("F18").Value = 427
